# bugfix for search in AddVendor
#
# - Remove the leftover test-customer rows (4 & 5) from the "clients_structured"
#   sheet, clearing their cell contents (this also drops the now-unused
#   "testcustomer1" / "testcustomer2" / "testcustomer1@gmail.com" shared strings).
# - Insert a blank row into the "vendors" lookup sheet (used by the AddVendor
#   search) at row 8, shifting the existing rows down by one.
# - Leave the UI pointed at "clients_structured" (instead of "SA-Technology")
#   as the active tab/sheet, with updated selections on the touched sheets.

$wb = $excel.ActiveWorkbook

# --- clients_structured: wipe the stray test-customer rows (11 & 12) -------
$wsClientsStructured = $wb.Worksheets.Item("clients_structured")

$wsClientsStructured.Range("A11").ClearContents()
$wsClientsStructured.Range("B11").Value = ""
$wsClientsStructured.Range("C11").Value = ""
$wsClientsStructured.Range("D11").Value = ""
$wsClientsStructured.Range("E11").ClearContents()

$wsClientsStructured.Range("A12").ClearContents()
$wsClientsStructured.Range("B12").Value = ""
$wsClientsStructured.Range("C12").Value = ""
$wsClientsStructured.Range("D12").Value = ""
$wsClientsStructured.Range("E12").ClearContents()

# --- vendors: insert a new blank row above row 8 (search range fix) --------
$wsVendors = $wb.Worksheets.Item("vendors")
$wsVendors.Rows.Item(8).Insert()
$wsVendors.Range("B14").Select()

# --- SA-Technology: move the selection, it's no longer the active tab ------
$wsSaTechnology = $wb.Worksheets.Item("SA-Technology")
$wsSaTechnology.Range("A7").Select()

# --- clients_structured becomes the active tab/sheet ------------------------
$wsClientsStructured.Activate()
